$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Structural changes ----

# Drop the extra 3rd data row (was "test" / "test 2" row).
$ws.Rows.Item(3).Delete()

# Insert a new first column for the Email field (shifts old A:R -> B:S).
$ws.Columns.Item(1).Insert()

# Insert a new column ahead of "Comments" (now column O after the shift above)
# for the new ratings-count field (shifts old N(Comments):R -> O:S further right).
$ws.Columns.Item(15).Insert()

# The two inserts above leave four trailing "timestamp" columns (Q,R,S,T) where
# the target layout only keeps two (Q,R) populated with fresh values. Remove the
# two surplus trailing columns.
$ws.Range("S1:T1").EntireColumn.Delete()

# ---- Formatting for the two newly-inserted header cells ----
# Reuse the existing header style (bold/border/centered) by copying format
# from a neighboring header cell instead of re-creating it.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---- Row 1: headers ----
$ws.Cells.Item(1,1).Value  = "Email"
$ws.Cells.Item(1,2).Value  = "SME"
$ws.Cells.Item(1,3).Value  = "Batch Name"
$ws.Cells.Item(1,4).Value  = "Course Event"
$ws.Cells.Item(1,5).Value  = "Camera On While Delivering"
$ws.Cells.Item(1,6).Value  = "Class Started on Time"
$ws.Cells.Item(1,7).Value  = "Zoom Poll Taken / Feedback Poll Taken"
$ws.Cells.Item(1,8).Value  = "Total attendees (online + offline)"
$ws.Cells.Item(1,9).Value  = "Resolution of Non Tech query"
$ws.Cells.Item(1,10).Value = "Resolution of Tech query"
$ws.Cells.Item(1,11).Value = "Refer and earn slide shown"
$ws.Cells.Item(1,12).Value = "Participant Engagement"
$ws.Cells.Item(1,13).Value = "Technical glitch (if any)"
$ws.Cells.Item(1,14).Value = "Was there any disruption during the session?"
$ws.Cells.Item(1,15).Value = "How many ratings less than 4 for today's session? (in any category)"
$ws.Cells.Item(1,16).Value = "Comments"
$ws.Cells.Item(1,17).Value = "2025-05-21 10:29:35"
$ws.Cells.Item(1,18).Value = "2025-05-21 10:33:14"

# ---- Row 2: data ----
$ws.Cells.Item(2,1).Value  = "shahul.s@skillfloor.com"
$ws.Cells.Item(2,2).Value  = "Shahul Hameed"
$ws.Cells.Item(2,3).Value  = "17-MAR-25-CDE-BUN-031-WDE2030-ONL31-MAR-25-CDE-BUN-031-WDE2030-ONL21-APR-25-CDE-BUN-031-WDE2030-ONL"
$ws.Cells.Item(2,4).Value  = "28-APR-25-DEA-135-WDE20-ONL (CONTINUE)"
$ws.Cells.Item(2,5).Value  = "Yes"
$ws.Cells.Item(2,6).Value  = "Yes"
$ws.Cells.Item(2,7).Value  = "Yes"
$ws.Cells.Item(2,8).Value  = 10
$ws.Cells.Item(2,9).Value  = "Yes"
$ws.Cells.Item(2,10).Value = "Yes"
$ws.Cells.Item(2,11).Value = "Yes"
$ws.Cells.Item(2,12).Value = "Yes"
$ws.Cells.Item(2,13).Value = "No"
$ws.Cells.Item(2,14).Value = "No"
$ws.Cells.Item(2,15).Value = 1
$ws.Cells.Item(2,16).Value = "Hey"
$ws.Cells.Item(2,17).Value = "Submitted"
$ws.Cells.Item(2,18).Value = "Submitted"
